$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("K-SVMeans")

# --- New "F" column (K=3, Aroon=5, Volume) for the existing K=2 block ---
$ws.Range("F4").Value = 59.42
$ws.Range("F5").Value = 55.67
$ws.Range("F6").Value = 51.87
$ws.Range("F7").Value = 57.61
$ws.Range("F8").Value = 55.84
$ws.Range("F9").Formula = "=AVERAGE(F4:F8)"

# Match the existing "best value" highlight style (reuse E8's cell format)
$ws.Range("E8").Copy()
$ws.Range("F6").PasteSpecial(-4122)

# Match the existing "Total" row style (reuse E9's cell format)
$ws.Range("E9").Copy()
$ws.Range("F9").PasteSpecial(-4122)

# --- New section: "K= 3, Aroon = 5, No volume" ---
$ws.Range("A11").Value = "K= 3, Aroon = 5, No volume"

$ws.Range("B12").Value = "period = 1"
$ws.Range("C12").Value = "period = 5"

$ws.Range("A13").Value = "BT6"
$ws.Range("B13").Value = 64.52

$ws.Range("A14").Value = "DHG"
$ws.Range("B14").Value = 59.14

$ws.Range("A15").Value = "FPT"
$ws.Range("B15").Value = 65.95

$ws.Range("A16").Value = "VIS"
$ws.Range("B16").Value = 51.89

$ws.Range("A17").Value = "VNM"
$ws.Range("B17").Value = 56.89

$ws.Range("A18").Value = "Total"
$ws.Range("B18").Formula = "=AVERAGE(B13:B17)"

# Highlight style on the new section's low values (reuse E8's cell format again)
$ws.Range("E8").Copy()
$ws.Range("B14:B15").PasteSpecial(-4122)

# "Total" row style (reuse E9's cell format again)
$ws.Range("E9").Copy()
$ws.Range("B18").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Page setup / view state ---
$ws.PageSetup.Orientation = 1

$ws.Activate()
$ws.Range("D15").Select()
